$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "758 nguyen binh, Hồ Chí Minh, Huyện Nhà Bè, Xã Nhơn Đức , Huyện Nhà Bè",
    "137A Nguyễn Văn Cừ - Long Biên - Hà Nội, Phường Ngọc Lâm, Quận Long Biên, Hà Nội",
    "Lữ Đoàn 45 , Cổ Đông , Sơn Tây, Hà nội, Hà Nội, Thành Phố Sơn Tây, Xã Cổ Đông , Thành Phố Sơn Tây",
    "Số 46 tổ 12, Phường Sài Đồng, Quận Long Biên, Hà Nội",
    "số 8 ngõ 191 phạm văn đồng, Phường Xuân Đỉnh, Quận Bắc Từ Liêm, Hà Nội",
    "60/1, Xã Bà Điểm, Huyện Hóc Môn, Hồ Chí Minh",
    "Cầu cảng số 12, đường hồ học lãm, p16, q8, tp.hcm, Phường 16, Quận 8, Hồ Chí Minh",
    "Số 4 tổ 51 cụm 8 Phú thượng Tây Hồ ( Chân cầu Nhật Tân ), Phường Phú Thượng, Quận Tây Hồ, Hà Nội",
    "so12 duong ht25 phuong hiep thanh q12.cty hoang vinh ha, Phường Hiệp Thành, Quận 12, Hồ Chí Minh",
    "Số 17 vạn phúc, Phường Kim Mã, Quận Ba Đình, Hà Nội",
    "Số 16 phố viên, gần viện bảo vệ Thực vật, Phường Cổ Nhuế 2, Quận Bắc Từ Liêm, Hà Nội",
    "68 Đường 1a Tổ3 kp Giản Dân phường Long Thạnh Mỹ Q9, Phường Long Thạnh Mỹ, Quận 9, Hồ Chí Minh"
)

# Replace the 9 old lookup results with the 12 new ones (rows 1-9 overwritten,
# rows 10-12 newly filled in).
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Append a new, identically-formatted blank row at the bottom of the results list.
$ws.Range("A99").Copy($ws.Range("A100")) | Out-Null
$ws.Range("A100").ClearContents() | Out-Null
$ws.Rows.Item(100).RowHeight = $ws.Rows.Item(99).RowHeight

$ws.Range("E100").Select() | Out-Null
